$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "current filter" -> "breadcrumb" (column A only; C8 already says breadcrumb) ---
$ws.Range("A8").Value = "breadcrumb"

# --- Row 9: same text ("collor filters"), but now bold on top of the existing yellow fill ---
$ws.Range("A9").Font.Bold = $true

# --- Rows 10-15: strip the "select " prefix from column A, and apply a bold
#     accent1 (theme4) colored font to match the new "collor filters" styling ---
$colorRows = 10,11,12,13,14,15
$colorNames = 'Black','Black','Gold','Pink','Grey','White'
for ($i = 0; $i -lt $colorRows.Length; $i++) {
    $r = $colorRows[$i]
    $ws.Cells.Item($r, 1).Value = $colorNames[$i]
    $ws.Cells.Item($r, 1).Font.Bold = $true
    $ws.Cells.Item($r, 1).Font.ThemeColor = 5
}

# --- Row 16 (new): "price range filters" section header, bold on yellow fill (like row 9) ---
$ws.Range("A16").Value = "price range filters"
$ws.Range("A16").Font.Bold = $true
$ws.Range("A16").Interior.Color = 65535

# --- Rows 17-22 (new): price range filter rows ---
$priceLabels = 'under 20','20-39','40-59','60-79','80-99','100-249'
$priceXpaths = @(
    "//div[contains(@class,'accordion-content is-open')]//span[contains(text(),'Under " + '$20' + "')]",
    "//div[contains(@class,'accordion-content is-open')]//span[contains(text(),'" + '$20 - $39' + "')]",
    "//div[contains(@class,'accordion-content is-open')]//span[contains(text(),'" + '$40 - $59' + "')]",
    "//div[contains(@class,'accordion-content is-open')]//span[contains(text(),'" + '$60 - $79' + "')]",
    "//div[contains(@class,'accordion-content is-open')]//span[contains(text(),'" + '$80 - $99' + "')]",
    "//div[contains(@class,'accordion-content is-open')]//span[contains(text(),'" + '$100 - $249' + "')]"
)
$priceExpected = 'under $20','$20-$39','$40-$59','$60-$79','$80-$99','$100-$249'

for ($i = 0; $i -lt 6; $i++) {
    $r = 17 + $i
    $ws.Cells.Item($r, 1).Value = $priceLabels[$i]
    $ws.Cells.Item($r, 2).Value = "xpath"
    $ws.Cells.Item($r, 3).Value = $priceXpaths[$i]
    $ws.Cells.Item($r, 4).Value = $priceExpected[$i]

    $ws.Cells.Item($r, 1).Font.Bold = $true
    $ws.Cells.Item($r, 1).Font.ThemeColor = 10
    $ws.Cells.Item($r, 4).Font.Bold = $true
    $ws.Cells.Item($r, 4).Font.ThemeColor = 10
}

# --- Column C width widened to fit the longer new xpath strings ---
$ws.Columns("C").ColumnWidth = 80.1

# --- Selection left on A8, matching the saved workbook state ---
$ws.Range("A8").Select() | Out-Null

Write-Output "done"
